{"js": "/*\n * Replace each arithmetic equation's text in the table with its updated\n * value, in document order. Every old/new value below appears exactly once\n * in the document, so a literal (non-wildcard) search-and-replace for each\n * pair is unambiguous and safe regardless of execution order.\n */\nconst pairs = [\n  [\"58+19=77\", \"34+27=61\"],\n  [\"12+67=79\", \"37-36=1\"],\n  [\"48-10=38\", \"0+41=41\"],\n  [\"44-42=2\", \"13+5=18\"],\n  [\"30-11=19\", \"35+29=64\"],\n  [\"99-85=14\", \"18-15=3\"],\n  [\"15+64=79\", \"15+6=21\"],\n  [\"97-40=57\", \"75-12=63\"],\n  [\"50+8=58\", \"29+55=84\"],\n  [\"25-5=20\", \"11+83=94\"],\n  [\"73-23=50\", \"57+17=74\"],\n  [\"24+31=55\", \"91-87=4\"],\n  [\"63-45=18\", \"86+5=91\"],\n  [\"4+39=43\", \"38+40=78\"],\n  [\"10-1=9\", \"45+33=78\"],\n  [\"39+22=61\", \"36-28=8\"],\n  [\"6+8=14\", \"33-29=4\"],\n  [\"38+41=79\", \"41-34=7\"],\n  [\"84-59=25\", \"2+92=94\"],\n  [\"97-26=71\", \"3+0=3\"],\n  [\"40+6=46\", \"98-36=62\"],\n  [\"70+12=82\", \"93-41=52\"],\n  [\"53+0=53\", \"48-18=30\"],\n  [\"61-36=25\", \"89+8=97\"],\n  [\"40+44=84\", \"89-16=73\"],\n  [\"31-5=26\", \"10+25=35\"],\n  [\"89-80=9\", \"16+51=67\"],\n  [\"75-2=73\", \"60-44=16\"],\n  [\"20+66=86\", \"37+36=73\"],\n  [\"42-11=31\", \"55-38=17\"],\n  [\"22-20=2\", \"20+57=77\"],\n  [\"74-8=66\", \"65-43=22\"],\n  [\"50+40=90\", \"95-26=69\"],\n  [\"53+6=59\", \"12+38=50\"],\n  [\"61-49=12\", \"77-19=58\"],\n  [\"55-39=16\", \"3+43=46\"],\n  [\"49+13=62\", \"20+46=66\"],\n  [\"1+85=86\", \"44+48=92\"],\n  [\"92-91=1\", \"24+51=75\"],\n  [\"42-21=21\", \"15+47=62\"],\n  [\"46-29=17\", \"5-5=0\"],\n  [\"93-92=1\", \"42+2=44\"],\n  [\"45-29=16\", \"24+17=41\"],\n  [\"55+40=95\", \"16+0=16\"],\n  [\"0+46=46\", \"95-58=37\"],\n  [\"46+40=86\", \"34+45=79\"],\n  [\"2+11=13\", \"41-30=11\"],\n  [\"31-6=25\", \"96-59=37\"],\n  [\"76-19=57\", \"31-2=29\"],\n  [\"54+40=94\", \"89-24=65\"],\n  [\"79-36=43\", \"34-33=1\"],\n  [\"43-9=34\", \"20+31=51\"],\n  [\"41-26=15\", \"90-24=66\"],\n  [\"94-16=78\", \"63-24=39\"],\n  [\"26+52=78\", \"72+3=75\"],\n  [\"84-70=14\", \"70-9=61\"],\n  [\"42-19=23\", \"56+13=69\"],\n  [\"16+66=82\", \"34+59=93\"],\n  [\"28+69=97\", \"74-65=9\"],\n  [\"73-57=16\", \"0+31=31\"],\n  [\"55+29=84\", \"18+37=55\"],\n  [\"84-83=1\", \"68+27=95\"],\n  [\"58-46=12\", \"9+37=46\"],\n  [\"75+10=85\", \"14-13=1\"],\n  [\"10+8=18\", \"40-24=16\"],\n  [\"67-8=59\", \"69-5=64\"],\n  [\"7+66=73\", \"71-10=61\"],\n  [\"16+73=89\", \"38+28=66\"],\n  [\"31+6=37\", \"18+64=82\"],\n  [\"79+2=81\", \"62+23=85\"],\n  [\"27-3=24\", \"91-82=9\"],\n  [\"18+50=68\", \"59+38=97\"],\n  [\"69-32=37\", \"57+37=94\"],\n  [\"22+5=27\", \"65-60=5\"],\n  [\"37+58=95\", \"10+56=66\"],\n  [\"5+74=79\", \"45-9=36\"],\n  [\"67+11=78\", \"35+13=48\"],\n  [\"1+37=38\", \"63-25=38\"],\n  [\"6+57=63\", \"88-44=44\"],\n  [\"84-58=26\", \"66-7=59\"],\n  [\"40+19=59\", \"9+87=96\"],\n  [\"25+25=50\", \"5+20=25\"],\n  [\"88-18=70\", \"11+34=45\"],\n  [\"73-28=45\", \"27-10=17\"],\n  [\"78-49=29\", \"35-13=22\"],\n  [\"2+97=99\", \"67-29=38\"],\n  [\"14+74=88\", \"38+54=92\"],\n  [\"24+48=72\", \"70-26=44\"],\n  [\"77-14=63\", \"31-4=27\"],\n  [\"13+65=78\", \"7+59=66\"],\n  [\"73-33=40\", \"32-31=1\"],\n  [\"18-14=4\", \"50+18=68\"],\n  [\"84-2=82\", \"78+12=90\"],\n  [\"50+21=71\", \"2+16=18\"],\n  [\"71-35=36\", \"82-72=10\"],\n  [\"53-28=25\", \"35+35=70\"],\n  [\"46+10=56\", \"76-65=11\"],\n  [\"41-23=18\", \"86+8=94\"],\n  [\"87-8=79\", \"15+0=15\"],\n  [\"64+2=66\", \"97-45=52\"]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true\n  });\n  results.load(\"items,text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n\n  // Replace the matched range's text while keeping its run formatting\n  // (font, size, etc.) untouched.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace each arithmetic equation's text in the table with its updated\n# value, in document order. Every old/new value below appears exactly once\n# in the document, so a literal (non-wildcard) find-and-replace for each\n# pair is unambiguous and safe regardless of execution order.\n$pairs = @(\n    @(\"58+19=77\", \"34+27=61\"),\n    @(\"12+67=79\", \"37-36=1\"),\n    @(\"48-10=38\", \"0+41=41\"),\n    @(\"44-42=2\", \"13+5=18\"),\n    @(\"30-11=19\", \"35+29=64\"),\n    @(\"99-85=14\", \"18-15=3\"),\n    @(\"15+64=79\", \"15+6=21\"),\n    @(\"97-40=57\", \"75-12=63\"),\n    @(\"50+8=58\", \"29+55=84\"),\n    @(\"25-5=20\", \"11+83=94\"),\n    @(\"73-23=50\", \"57+17=74\"),\n    @(\"24+31=55\", \"91-87=4\"),\n    @(\"63-45=18\", \"86+5=91\"),\n    @(\"4+39=43\", \"38+40=78\"),\n    @(\"10-1=9\", \"45+33=78\"),\n    @(\"39+22=61\", \"36-28=8\"),\n    @(\"6+8=14\", \"33-29=4\"),\n    @(\"38+41=79\", \"41-34=7\"),\n    @(\"84-59=25\", \"2+92=94\"),\n    @(\"97-26=71\", \"3+0=3\"),\n    @(\"40+6=46\", \"98-36=62\"),\n    @(\"70+12=82\", \"93-41=52\"),\n    @(\"53+0=53\", \"48-18=30\"),\n    @(\"61-36=25\", \"89+8=97\"),\n    @(\"40+44=84\", \"89-16=73\"),\n    @(\"31-5=26\", \"10+25=35\"),\n    @(\"89-80=9\", \"16+51=67\"),\n    @(\"75-2=73\", \"60-44=16\"),\n    @(\"20+66=86\", \"37+36=73\"),\n    @(\"42-11=31\", \"55-38=17\"),\n    @(\"22-20=2\", \"20+57=77\"),\n    @(\"74-8=66\", \"65-43=22\"),\n    @(\"50+40=90\", \"95-26=69\"),\n    @(\"53+6=59\", \"12+38=50\"),\n    @(\"61-49=12\", \"77-19=58\"),\n    @(\"55-39=16\", \"3+43=46\"),\n    @(\"49+13=62\", \"20+46=66\"),\n    @(\"1+85=86\", \"44+48=92\"),\n    @(\"92-91=1\", \"24+51=75\"),\n    @(\"42-21=21\", \"15+47=62\"),\n    @(\"46-29=17\", \"5-5=0\"),\n    @(\"93-92=1\", \"42+2=44\"),\n    @(\"45-29=16\", \"24+17=41\"),\n    @(\"55+40=95\", \"16+0=16\"),\n    @(\"0+46=46\", \"95-58=37\"),\n    @(\"46+40=86\", \"34+45=79\"),\n    @(\"2+11=13\", \"41-30=11\"),\n    @(\"31-6=25\", \"96-59=37\"),\n    @(\"76-19=57\", \"31-2=29\"),\n    @(\"54+40=94\", \"89-24=65\"),\n    @(\"79-36=43\", \"34-33=1\"),\n    @(\"43-9=34\", \"20+31=51\"),\n    @(\"41-26=15\", \"90-24=66\"),\n    @(\"94-16=78\", \"63-24=39\"),\n    @(\"26+52=78\", \"72+3=75\"),\n    @(\"84-70=14\", \"70-9=61\"),\n    @(\"42-19=23\", \"56+13=69\"),\n    @(\"16+66=82\", \"34+59=93\"),\n    @(\"28+69=97\", \"74-65=9\"),\n    @(\"73-57=16\", \"0+31=31\"),\n    @(\"55+29=84\", \"18+37=55\"),\n    @(\"84-83=1\", \"68+27=95\"),\n    @(\"58-46=12\", \"9+37=46\"),\n    @(\"75+10=85\", \"14-13=1\"),\n    @(\"10+8=18\", \"40-24=16\"),\n    @(\"67-8=59\", \"69-5=64\"),\n    @(\"7+66=73\", \"71-10=61\"),\n    @(\"16+73=89\", \"38+28=66\"),\n    @(\"31+6=37\", \"18+64=82\"),\n    @(\"79+2=81\", \"62+23=85\"),\n    @(\"27-3=24\", \"91-82=9\"),\n    @(\"18+50=68\", \"59+38=97\"),\n    @(\"69-32=37\", \"57+37=94\"),\n    @(\"22+5=27\", \"65-60=5\"),\n    @(\"37+58=95\", \"10+56=66\"),\n    @(\"5+74=79\", \"45-9=36\"),\n    @(\"67+11=78\", \"35+13=48\"),\n    @(\"1+37=38\", \"63-25=38\"),\n    @(\"6+57=63\", \"88-44=44\"),\n    @(\"84-58=26\", \"66-7=59\"),\n    @(\"40+19=59\", \"9+87=96\"),\n    @(\"25+25=50\", \"5+20=25\"),\n    @(\"88-18=70\", \"11+34=45\"),\n    @(\"73-28=45\", \"27-10=17\"),\n    @(\"78-49=29\", \"35-13=22\"),\n    @(\"2+97=99\", \"67-29=38\"),\n    @(\"14+74=88\", \"38+54=92\"),\n    @(\"24+48=72\", \"70-26=44\"),\n    @(\"77-14=63\", \"31-4=27\"),\n    @(\"13+65=78\", \"7+59=66\"),\n    @(\"73-33=40\", \"32-31=1\"),\n    @(\"18-14=4\", \"50+18=68\"),\n    @(\"84-2=82\", \"78+12=90\"),\n    @(\"50+21=71\", \"2+16=18\"),\n    @(\"71-35=36\", \"82-72=10\"),\n    @(\"53-28=25\", \"35+35=70\"),\n    @(\"46+10=56\", \"76-65=11\"),\n    @(\"41-23=18\", \"86+8=94\"),\n    @(\"87-8=79\", \"15+0=15\"),\n    @(\"64+2=66\", \"97-45=52\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1              # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    # 2 = wdReplaceAll (only one occurrence exists, so this replaces it once)\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
